$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.476.85"
$ws.Range("E2").Value = "'  +0.37%  "
$ws.Range("D3").Value = "'1.852.56"
$ws.Range("E3").Value = "'  +0.41%  "
$ws.Range("E5").Value = "'  +0.88%  "
$ws.Range("D6").Value = "'0.6304"
$ws.Range("E6").Value = "'  +0.24%  "
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("D8").Value = "'0.07676"
$ws.Range("E8").Value = "'  +1.72%  "
$ws.Range("D9").Value = "'0.2936"
$ws.Range("E9").Value = "'  -0.27%  "
$ws.Range("E10").Value = "'  +0.52%  "
$ws.Range("D11").Value = "'0.07752"
$ws.Range("E11").Value = "'  +0.74%  "
$ws.Range("D12").Value = "'1.871.29"
$ws.Range("E12").Value = "'  +1.44%  "
$ws.Range("E13").Value = "'  +1.20%  "
$ws.Range("D14").Value = "'0.6814"
$ws.Range("E14").Value = "'  +0.53%  "
$ws.Range("D15").Value = "'0.00001067"
$ws.Range("E15").Value = "'  +4.65%  "
$ws.Range("E16").Value = "'  +0.88%  "
$ws.Range("D17").Value = "'2.117.37"
$ws.Range("E17").Value = "'  +1.26%  "
$ws.Range("D18").Value = "'6.179"
$ws.Range("E18").Value = "'  +0.84%  "
$ws.Range("D19").Value = "'29.505.78"
$ws.Range("E19").Value = "'  +0.35%  "
$ws.Range("D20").Value = "'229.59"
$ws.Range("E20").Value = "'  +0.62%  "
$ws.Range("E21").Value = "'  +0.46%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "'  +0.08%  "
$ws.Range("E23").Value = "'  +0.16%  "
$ws.Range("E24").Value = "'  +0.06%  "
$ws.Range("D25").Value = "'157.00"
$ws.Range("E25").Value = "'  +0.20%  "
$ws.Range("D26").Value = "'0.1383"
$ws.Range("E26").Value = "'  -0.63%  "
$ws.Range("D27").Value = "'8.404"
$ws.Range("E27").Value = "'  +0.65%  "
$ws.Range("E28").Value = "'  +0.51%  "
$ws.Range("D29").Value = "'1.332"
$ws.Range("E29").Value = "'  +5.11%  "
$ws.Range("D30").Value = "'1.467"
$ws.Range("E30").Value = "'  +0.73%  "
$ws.Range("D31").Value = "'0.05690"
$ws.Range("E31").Value = "'  +1.15%  "
$ws.Range("D32").Value = "'4.138"
$ws.Range("E32").Value = "'  +0.49%  "
$ws.Range("E33").Value = "'  +0.20%  "
$ws.Range("E34").Value = "'  +0.99%  "
$ws.Range("D35").Value = "'1.166"
$ws.Range("E35").Value = "'  +1.04%  "
$ws.Range("D36").Value = "'0.7093"
$ws.Range("E36").Value = "'  -0.89%  "
$ws.Range("D37").Value = "'2.587"
$ws.Range("E37").Value = "'  -0.30%  "
$ws.Range("D38").Value = "'2.782"
$ws.Range("E38").Value = "'  +0.40%  "
$ws.Range("E39").Value = "'  -0.70%  "
$ws.Range("D40").Value = "'1.220.79"
$ws.Range("D41").Value = "'6.558"
$ws.Range("E41").Value = "'  +5.50%  "
$ws.Range("D42").Value = "'0.9084"
$ws.Range("E42").Value = "'  +0.78%  "
$ws.Range("E43").Value = "'  +0.11%  "
$ws.Range("B44").Value = "'Quant"
$ws.Range("C44").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'101.64"
$ws.Range("E44").Value = "'  -0.15%  "
$ws.Range("B45").Value = "'Aave"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'66.52"
$ws.Range("E45").Value = "'  +1.01%  "
$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000120"
$ws.Range("E46").Value = "'  +0.87%  "
$ws.Range("B47").Value = "'Aptos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.137"
$ws.Range("E47").Value = "'  +0.43%  "
$ws.Range("B48").Value = "'TheSandbox"
$ws.Range("C48").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "'0.4025"
$ws.Range("E48").Value = "'  +0.73%  "
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.049"
$ws.Range("E49").Value = "'  +1.09%  "
$ws.Range("B50").Value = "'RenderToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.688"
$ws.Range("E50").Value = "'  +0.41%  "
$ws.Range("B51").Value = "'Algorand"
$ws.Range("C51").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1146"
$ws.Range("E51").Value = "'  +2.58%  "
